$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark from its old location
#    (right before "The direct..."). We'll recreate it at the end
#    of the new paragraph we are about to add.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Split the existing paragraph into two: insert a paragraph break
#    right after "...generate the desire product object." and add
#    the new sentences as a new paragraph.
# ------------------------------------------------------------------
$searchRange = $d.Range(0, $d.Content.End)
$found = $searchRange.Find.Execute("desire product object.")
if (-not $found) {
    throw "Could not find anchor text 'desire product object.'"
}
$endOfOldParagraph = $searchRange.End
$searchRange.Collapse(0)              # wdCollapseEnd
$searchRange.InsertParagraphAfter()

$newParaRange = $d.Range($endOfOldParagraph + 1, $endOfOldParagraph + 1)
$newParaText = "Our implementation of the builder pattern uses a Row Builder to identify the attributes relevant to the requested type of file. An Attribute Receiver then returns the values associated with each element."
$newParaRange.InsertAfter($newParaText)

# ------------------------------------------------------------------
# 3. Re-add the "_GoBack" bookmark, collapsed, at the very end of the
#    document (end of the new paragraph). Adding a bookmark exactly
#    at the current end-of-story position is unreliable in this
#    runtime, so we temporarily pad the document, add the bookmark,
#    then strip the padding back out again.
# ------------------------------------------------------------------
$marker = [char]0x2060   # rare, invisible "word joiner" character
$padding = "$marker$marker"

$bookmarkPos = $d.Content.End
$tail = $d.Content
$tail.InsertAfter($padding)

$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$paddedEnd = $d.Content.End
$padRange = $d.Range($bookmarkPos, $paddedEnd)
$padRange.Delete()

$cleanupRange = $d.Range(0, $d.Content.End)
$cleanupRange.Find.Execute([string]$marker, $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

Write-Host "done"
